$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 45 (Arveja Verde, Perfection,
# Provincia de Huasco). All subsequent rows (old 45..127) shift down one row
# (to 46..128), exactly as Excel does with a real row insert.
$ws.Rows("45:45").Insert()

$ws.Cells.Item(45,1).Value  = 4
$ws.Cells.Item(45,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45,3).Value  = "Los Lagos"
$ws.Cells.Item(45,4).Value  = 44771
$ws.Cells.Item(45,5).Value  = 10
$ws.Cells.Item(45,6).Value  = 100112022
$ws.Cells.Item(45,7).Value  = "Arveja Verde"
$ws.Cells.Item(45,8).Value  = "Perfection"
$ws.Cells.Item(45,9).Value  = "Primera"
$ws.Cells.Item(45,10).Value = 80
$ws.Cells.Item(45,11).Value = 35000
$ws.Cells.Item(45,12).Value = 35000
$ws.Cells.Item(45,13).Value = 35000
$ws.Cells.Item(45,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(45,15).Value = "Provincia de Huasco"
$ws.Cells.Item(45,16).Value = 1400
$ws.Cells.Item(45,17).Value = 25
$ws.Cells.Item(45,18).Value = "Hortaliza"
